$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N19").Value = -3294.7
$ws.Range("J19").Value = 2944.7
$ws.Range("H19").Value = 3497.7083
$ws.Range("L19").Value = 2944.7
$ws.Range("J51").Value = 2361
$ws.Range("H51").Value = 2361
$ws.Range("N51").Value = -3329
$ws.Range("L51").Value = 2361
$ws.Range("M76").Value = -3118
$ws.Range("H76").Value = 6476.4443
$ws.Range("K76").Value = 3433
$ws.Range("J76").Value = 7998.1665
$ws.Range("I76").Value = 3433
$ws.Range("N76").Value = -8628.166499999999
$ws.Range("L76").Value = 7998.1665
$ws.Range("H79").Value = 6476.4443
$ws.Range("M79").Value = -2341
$ws.Range("I79").Value = 3433
$ws.Range("L79").Value = 7998.1665
$ws.Range("K79").Value = 3433
$ws.Range("N79").Value = -10182.1665
$ws.Range("J79").Value = 7998.1665
$ws.Range("H86").Value = 3819.7058
$ws.Range("I86").Value = 1619.5
$ws.Range("M86").Value = -496.5
$ws.Range("K86").Value = 1619.5
$ws.Range("L87").Value = 99999
$ws.Range("J87").Value = 99999
$ws.Range("N87").Value = -102495
$ws.Range("H87").Value = 99999
$ws.Range("M89").Value = -2481.5
$ws.Range("H89").Value = 3819.7058
$ws.Range("K89").Value = 8097.5
$ws.Range("I89").Value = 1619.5
$ws.Range("L90").Value = 299997
$ws.Range("N90").Value = -312477
$ws.Range("H90").Value = 99999
$ws.Range("J90").Value = 99999
$ws.Range("N100").Value = -9514.1
$ws.Range("J100").Value = 8432.1
$ws.Range("L100").Value = 8432.1
$ws.Range("H100").Value = 6502.8823
$ws.Range("H132").Value = 2621.0417
$ws.Range("I132").Value = 1810.9697
$ws.Range("K132").Value = 5432.909100000001
$ws.Range("M132").Value = -2902.909100000001
$ws.Range("M137").Value = -272726430
$ws.Range("I137").Value = 90909660
$ws.Range("H137").Value = 47620870
$ws.Range("K137").Value = 272728980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 1151.62
$ws.Range("H32").Value = 1134.4423
$ws.Range("M32").Value = -864.6199999999999
$ws.Range("I32").Value = 1151.62
$ws.Range("J88").Value = 1224.2858
$ws.Range("N88").Value = -2036.2858
$ws.Range("L88").Value = 1224.2858
$ws.Range("H88").Value = 1003.26666
$ws.Range("N91").Value = -4032.2858
$ws.Range("J91").Value = 1224.2858
$ws.Range("L91").Value = 1224.2858
$ws.Range("H91").Value = 1003.26666
$ws.Range("J112").Value = 15000
$ws.Range("N112").Value = -17954
$ws.Range("H112").Value = 15000
$ws.Range("L112").Value = 15000
$ws.Range("M122").Value = -2267.5
$ws.Range("K122").Value = 4717.5
$ws.Range("H122").Value = 2227
$ws.Range("I122").Value = 1572.5
$ws.Range("H132").Value = 7946.7144
$ws.Range("N132").Value = -33188
$ws.Range("J132").Value = 9376
$ws.Range("I132").Value = 6874.75
$ws.Range("L132").Value = 28128
$ws.Range("K132").Value = 20624.25
$ws.Range("M132").Value = -18094.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1534.6666
$ws.Range("I86").Value = 1487.6666
$ws.Range("N86").Value = -3837.0667
$ws.Range("M86").Value = -364.6666
$ws.Range("L86").Value = 1591.0667
$ws.Range("J86").Value = 1591.0667
$ws.Range("K86").Value = 1487.6666
$ws.Range("J89").Value = 1591.0667
$ws.Range("M89").Value = -1822.333000000001
$ws.Range("L89").Value = 7955.333500000001
$ws.Range("H89").Value = 1534.6666
$ws.Range("K89").Value = 7438.333000000001
$ws.Range("N89").Value = -19187.3335
$ws.Range("I89").Value = 1487.6666
$ws.Range("H105").Value = 1198.5
$ws.Range("M105").Value = 848.6667
$ws.Range("L105").Value = 1498.6666
$ws.Range("K105").Value = 898.3333
$ws.Range("N105").Value = -4992.6666
$ws.Range("J105").Value = 1498.6666
$ws.Range("I105").Value = 898.3333
$ws.Range("M134").Value = -894.2552999999998
$ws.Range("K134").Value = 3429.2553
$ws.Range("I134").Value = 1143.0851
$ws.Range("H134").Value = 1856.4259

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K94").Value = 981
$ws.Range("J94").Value = 1482.2222
$ws.Range("M94").Value = -530
$ws.Range("H94").Value = 1303.2142
$ws.Range("I94").Value = 981
$ws.Range("L94").Value = 1482.2222
$ws.Range("N94").Value = -2384.2222
$ws.Range("I99").Value = 1670.375
$ws.Range("K99").Value = 1670.375
$ws.Range("M99").Value = -172.375
$ws.Range("H99").Value = 2131.1667
$ws.Range("H105").Value = 1294.9
$ws.Range("M105").Value = 982.5714
$ws.Range("L105").Value = 2532.6667
$ws.Range("K105").Value = 764.4286
$ws.Range("N105").Value = -6026.6667
$ws.Range("J105").Value = 2532.6667
$ws.Range("I105").Value = 764.4286
$ws.Range("I126").Value = 1670.375
$ws.Range("H126").Value = 2131.1667
$ws.Range("K126").Value = 5011.125
$ws.Range("M126").Value = -2541.125
$ws.Range("H132").Value = 133338850
$ws.Range("N132").Value = -27558.125
$ws.Range("J132").Value = 7499.375
$ws.Range("I132").Value = 285717540
$ws.Range("L132").Value = 22498.125
$ws.Range("K132").Value = 857152620
$ws.Range("M132").Value = -857150090
$ws.Range("J134").Value = 3325.25
$ws.Range("M134").Value = -4955.0625
$ws.Range("N134").Value = -15045.75
$ws.Range("L134").Value = 9975.75
$ws.Range("K134").Value = 7490.0625
$ws.Range("I134").Value = 2496.6875
$ws.Range("H134").Value = 2662.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 278913150
$ws.Range("N32").Value = -976195676
$ws.Range("J32").Value = 325398370
$ws.Range("L32").Value = 976195110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M97").Value = -376
$ws.Range("I97").Value = 872
$ws.Range("H97").Value = 929.45715
$ws.Range("K97").Value = 872
$ws.Range("N122").Value = -44907.001
$ws.Range("J122").Value = 13335.667
$ws.Range("M122").Value = -141989.641
$ws.Range("K122").Value = 144439.641
$ws.Range("L122").Value = 40007.001
$ws.Range("H122").Value = 43969.24
$ws.Range("I122").Value = 48146.547
$ws.Range("H132").Value = 13709963
$ws.Range("N132").Value = -43469.945
$ws.Range("J132").Value = 12803.315
$ws.Range("I132").Value = 18529334
$ws.Range("L132").Value = 38409.945
$ws.Range("K132").Value = 55588002
$ws.Range("M132").Value = -55585472

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 3491.7727
$ws.Range("K40").Value = 3491.7727
$ws.Range("M40").Value = -3355.7727
$ws.Range("H40").Value = 3335.3333
$ws.Range("M55").Value = -20.69999999999999
$ws.Range("K55").Value = 193.7
$ws.Range("I55").Value = 193.7
$ws.Range("H55").Value = 205.85
$ws.Range("M61").Value = -665.6667
$ws.Range("K61").Value = 867.6667
$ws.Range("J61").Value = 775.5
$ws.Range("H61").Value = 830.8
$ws.Range("I61").Value = 867.6667
$ws.Range("L61").Value = 775.5
$ws.Range("N61").Value = -1179.5
$ws.Range("K93").Value = 2300
$ws.Range("I93").Value = 2300
$ws.Range("M93").Value = -1052
$ws.Range("H93").Value = 2300
$ws.Range("J104").Value = 98481
$ws.Range("L104").Value = 98481
$ws.Range("N104").Value = -105469
$ws.Range("H104").Value = 98481
$ws.Range("L110").Value = 80644
$ws.Range("H110").Value = 80644
$ws.Range("N110").Value = -88824
$ws.Range("J110").Value = 80644
$ws.Range("L113").Value = 775.5
$ws.Range("I113").Value = 867.6667
$ws.Range("N113").Value = -5115.5
$ws.Range("H113").Value = 830.8
$ws.Range("K113").Value = 867.6667
$ws.Range("J113").Value = 775.5
$ws.Range("M113").Value = 1302.3333
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("N122").Value = -23897.5
$ws.Range("J122").Value = 6332.5
$ws.Range("M122").Value = -11206.3165
$ws.Range("K122").Value = 13656.3165
$ws.Range("L122").Value = 18997.5
$ws.Range("H122").Value = 4979.4
$ws.Range("I122").Value = 4552.1055
$ws.Range("N136").Value = -28396.2
$ws.Range("L136").Value = 23296.2
$ws.Range("H136").Value = 3566.2778
$ws.Range("I136").Value = 1951.2307
$ws.Range("M136").Value = -3303.6921
$ws.Range("K136").Value = 5853.6921
$ws.Range("J136").Value = 7765.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N81").Value = -21549.714
$ws.Range("L81").Value = 19427.714
$ws.Range("H81").Value = 10465.833
$ws.Range("J81").Value = 9713.857
$ws.Range("N84").Value = -107746.57
$ws.Range("J84").Value = 9713.857
$ws.Range("H84").Value = 10465.833
$ws.Range("L84").Value = 97138.57000000001
$ws.Range("L113").Value = 8246.143199999999
$ws.Range("I113").Value = 504.66666
$ws.Range("N113").Value = -12586.1432
$ws.Range("H113").Value = 1713
$ws.Range("K113").Value = 1513.99998
$ws.Range("J113").Value = 2748.7144
$ws.Range("M113").Value = 656.0000199999999
$ws.Range("J118").Value = 98000
$ws.Range("N118").Value = -101314
$ws.Range("L118").Value = 98000
$ws.Range("H118").Value = 98000
$ws.Range("N122").Value = -27392.5
$ws.Range("J122").Value = 7497.5
$ws.Range("M122").Value = -8497.193499999999
$ws.Range("K122").Value = 10947.1935
$ws.Range("L122").Value = 22492.5
$ws.Range("H122").Value = 3882.303
$ws.Range("I122").Value = 3649.0645
$ws.Range("H132").Value = 14292465
$ws.Range("N132").Value = -30564.875
$ws.Range("J132").Value = 8501.625
$ws.Range("I132").Value = 33337748
$ws.Range("L132").Value = 25504.875
$ws.Range("K132").Value = 100013244
$ws.Range("M132").Value = -100010714
$ws.Range("N135").Value = -130890
$ws.Range("J135").Value = 120750
$ws.Range("L135").Value = 120750
$ws.Range("H135").Value = 120750
$ws.Range("H136").Value = 15878335
$ws.Range("I136").Value = 18523570
$ws.Range("M136").Value = -55568160
$ws.Range("K136").Value = 55570710

Write-Output "Applied all Excalibur_Profits updates"